$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths for columns A and B from 66 to 81
# (ColumnWidth assignment adds a fixed 5/6 character padding offset when
# round-tripped to the OOXML "width" attribute, so subtract it here to
# land on exactly 81 in the saved file)
$ws.Columns.Item(1).ColumnWidth = 80.16666666666667
$ws.Columns.Item(2).ColumnWidth = 80.16666666666667

# Update text values in A2 and B2
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/scheduleTestRunAndConfigureEnvironment-test-data"
$ws.Range("B2").Value = "Data Files/AI-Generated/Common/scheduleTestRunAndConfigureEnvironment-test-data"
